$p = $ppt.ActivePresentation

# --- Slide 1 ("Lead Score Assignment" title slide): fix a name in the
# credits text box from "Shakti Tirni " to "Shaktiprasanna".
# The text box is shape 2 ("Text Placeholder 2"); the name is the 3rd
# bullet paragraph, originally split across three runs:
#   [1-7]  "Shakti "   (no err flag)
#   [8-12] "Tirni"     (err="1" - flagged by the spell checker)
#   [13]   " "         (trailing space, no err flag)
# The edit below removes the leading "Shakti " run, rewrites the
# "Tirni" run (keeping its own run properties, including err="1") to
# the full corrected name, and removes the trailing space run - giving
# a single run "Shaktiprasanna" with the err="1" flag preserved.
$slide1 = $p.Slides.Item(1)
$nameShape = $slide1.Shapes.Item(2)
$nameTextRange = $nameShape.TextFrame.TextRange
$namePara = $nameTextRange.Paragraphs(3)

$leadingRun = $namePara.Characters(1, 7)
$leadingRun.Text = ""

$middleRun = $namePara.Characters(1, 5)
$middleRun.Text = "Shaktiprasanna"

$trailingRun = $namePara.Characters(15, 1)
$trailingRun.Text = ""
